# "cierre 18 FEB 22" -- roll the payroll receipt ("RECIBO DE NOMINA") sheet
# forward one week: SEMANA 06 (7-13 FEB 2022) -> SEMANA 07 (14-20 FEB 2022),
# and update this period's hours/amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")
$ws.Activate()

# Header text used throughout the sheet (B9 is the master cell; H9, B27,
# H27, B43, H43 and B60 all carry formulas that point back to it, so they
# ripple automatically on recalculation).
$ws.Range("B9").Value = "SEMANA   07  DEL    14      Al   20   DE   FEBRERO          2022"

# Top block (days / extra-hours pay) for the first employee of the week.
$ws.Range("J3").Value = 2
$ws.Range("K3").Value = 867
$ws.Range("K4").Value = 0

# Second block: loan/advance deduction ("PRESTAMO") cleared to 0.
$ws.Range("E25").Value = 0

# Keep the view roughly where the edit happened.
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("D16:D17").Select()

$wb.Application.Calculate()
